$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing organizational form (D3) and name (E4) values
$ws.Range("D3").Value = "ПАО"
$ws.Range("E4").Value = "Имя 3"

# Update the active selection to match the saved view state
$ws.Range("E2:E4").Select()
